$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "model_22_3_24"
$ws.Cells.Item(2, 2).Value = [double]"0.9999725909228852"
$ws.Cells.Item(2, 3).Value = [double]"0.9990573586657896"
$ws.Cells.Item(2, 4).Value = [double]"0.9998944825218106"
$ws.Cells.Item(2, 5).Value = [double]"0.9999789093919152"
$ws.Cells.Item(2, 6).Value = [double]"0.9999182320163053"
$ws.Cells.Item(2, 7).Value = [double]"2.558517036724789e-05"
$ws.Cells.Item(2, 8).Value = [double]"0.0008799143083129368"
$ws.Cells.Item(2, 9).Value = [double]"0.0001472337679308854"
$ws.Cells.Item(2, 10).Value = [double]"1.15175906314148e-05"
$ws.Cells.Item(2, 11).Value = [double]"7.937567928115011e-05"
$ws.Cells.Item(2, 12).Value = [double]"0.0003185012387814086"
$ws.Cells.Item(2, 13).Value = [double]"0.005058178562214652"
$ws.Cells.Item(2, 14).Value = [double]"1.000038695167692"
$ws.Cells.Item(2, 15).Value = [double]"0.005273515695481226"
$ws.Cells.Item(2, 16).Value = [double]"103.1469953137142"
$ws.Cells.Item(2, 17).Value = [double]"153.1209041333104"

$ws.Cells.Item(3, 1).Value = "model_22_3_23"
$ws.Cells.Item(3, 2).Value = [double]"0.9999728855296113"
$ws.Cells.Item(3, 3).Value = [double]"0.9990564175101938"
$ws.Cells.Item(3, 4).Value = [double]"0.9998959206341197"
$ws.Cells.Item(3, 5).Value = [double]"0.9999788746249005"
$ws.Cells.Item(3, 6).Value = [double]"0.9999192558112617"
$ws.Cells.Item(3, 7).Value = [double]"2.531016792015338e-05"
$ws.Cells.Item(3, 8).Value = [double]"0.0008807928357497334"
$ws.Cells.Item(3, 9).Value = [double]"0.0001452270985371712"
$ws.Cells.Item(3, 10).Value = [double]"1.153657691391651e-05"
$ws.Cells.Item(3, 11).Value = [double]"7.838183772554385e-05"
$ws.Cells.Item(3, 12).Value = [double]"0.0003213210083516168"
$ws.Cells.Item(3, 13).Value = [double]"0.005030921180077598"
$ws.Cells.Item(3, 14).Value = [double]"1.000038279252313"
$ws.Cells.Item(3, 15).Value = [double]"0.005245097910156038"
$ws.Cells.Item(3, 16).Value = [double]"103.1686086977828"
$ws.Cells.Item(3, 17).Value = [double]"153.142517517379"

$ws.Cells.Item(4, 1).Value = "model_22_3_22"
$ws.Cells.Item(4, 2).Value = [double]"0.9999732056265778"
$ws.Cells.Item(4, 3).Value = [double]"0.9990553558295356"
$ws.Cells.Item(4, 4).Value = [double]"0.9998975157236147"
$ws.Cells.Item(4, 5).Value = [double]"0.9999788269800344"
$ws.Cells.Item(4, 6).Value = [double]"0.9999203888037597"
$ws.Cells.Item(4, 7).Value = [double]"2.501137145256621e-05"
$ws.Cells.Item(4, 8).Value = [double]"0.000881783868041752"
$ws.Cells.Item(4, 9).Value = [double]"0.0001430013910945274"
$ws.Cells.Item(4, 10).Value = [double]"1.156259579689735e-05"
$ws.Cells.Item(4, 11).Value = [double]"7.728199344571237e-05"
$ws.Cells.Item(4, 12).Value = [double]"0.0003244507699737741"
$ws.Cells.Item(4, 13).Value = [double]"0.005001137015976088"
$ws.Cells.Item(4, 14).Value = [double]"1.000037827350714"
$ws.Cells.Item(4, 15).Value = [double]"0.005214045772527006"
$ws.Cells.Item(4, 16).Value = [double]"103.19235995682"
$ws.Cells.Item(4, 17).Value = [double]"153.1662687764162"

$ws.Cells.Item(5, 1).Value = "model_22_3_21"
$ws.Cells.Item(5, 2).Value = [double]"0.9999735505297692"
$ws.Cells.Item(5, 3).Value = [double]"0.9990541704589598"
$ws.Cells.Item(5, 4).Value = [double]"0.9998992754132662"
$ws.Cells.Item(5, 5).Value = [double]"0.9999787657025544"
$ws.Cells.Item(5, 6).Value = [double]"0.9999216335051023"
$ws.Cells.Item(5, 7).Value = [double]"2.468941946284989e-05"
$ws.Cells.Item(5, 8).Value = [double]"0.0008828903594425073"
$ws.Cells.Item(5, 9).Value = [double]"0.0001405460088942092"
$ws.Cells.Item(5, 10).Value = [double]"1.159605945649702e-05"
$ws.Cells.Item(5, 11).Value = [double]"7.607370861210372e-05"
$ws.Cells.Item(5, 12).Value = [double]"0.0003278703040278962"
$ws.Cells.Item(5, 13).Value = [double]"0.004968844882148152"
$ws.Cells.Item(5, 14).Value = [double]"1.000037340428561"
$ws.Cells.Item(5, 15).Value = [double]"0.005180378895708083"
$ws.Cells.Item(5, 16).Value = [double]"103.2182715358628"
$ws.Cells.Item(5, 17).Value = [double]"153.1921803554591"

$ws.Cells.Item(6, 1).Value = "model_22_3_20"
$ws.Cells.Item(6, 2).Value = [double]"0.9999739234437186"
$ws.Cells.Item(6, 3).Value = [double]"0.9990528355944712"
$ws.Cells.Item(6, 4).Value = [double]"0.9999012225100296"
$ws.Cells.Item(6, 5).Value = [double]"0.9999786923580651"
$ws.Cells.Item(6, 6).Value = [double]"0.9999230150121922"
$ws.Cells.Item(6, 7).Value = [double]"2.434132066012588e-05"
$ws.Cells.Item(6, 8).Value = [double]"0.0008841363968487559"
$ws.Cells.Item(6, 9).Value = [double]"0.0001378291282605888"
$ws.Cells.Item(6, 10).Value = [double]"1.163611291531681e-05"
$ws.Cells.Item(6, 11).Value = [double]"7.473262058795281e-05"
$ws.Cells.Item(6, 12).Value = [double]"0.0003316237571951391"
$ws.Cells.Item(6, 13).Value = [double]"0.004933692396180155"
$ws.Cells.Item(6, 14).Value = [double]"1.000036813961809"
$ws.Cells.Item(6, 15).Value = [double]"0.00514372989563675"
$ws.Cells.Item(6, 16).Value = [double]"103.246670426077"
$ws.Cells.Item(6, 17).Value = [double]"153.2205792456732"

$ws.Cells.Item(7, 1).Value = "model_22_3_19"
$ws.Cells.Item(7, 2).Value = [double]"0.9999743149631809"
$ws.Cells.Item(7, 3).Value = [double]"0.9990513275006083"
$ws.Cells.Item(7, 4).Value = [double]"0.9999033445893806"
$ws.Cells.Item(7, 5).Value = [double]"0.9999785891037163"
$ws.Cells.Item(7, 6).Value = [double]"0.9999245111115425"
$ws.Cells.Item(7, 7).Value = [double]"2.397585442780653e-05"
$ws.Cells.Item(7, 8).Value = [double]"0.0008855441362721773"
$ws.Cells.Item(7, 9).Value = [double]"0.0001348680857483063"
$ws.Cells.Item(7, 10).Value = [double]"1.169250016196276e-05"
$ws.Cells.Item(7, 11).Value = [double]"7.328029295513451e-05"
$ws.Cells.Item(7, 12).Value = [double]"0.0003357708350012782"
$ws.Cells.Item(7, 13).Value = [double]"0.00489651451828814"
$ws.Cells.Item(7, 14).Value = [double]"1.000036261228451"
$ws.Cells.Item(7, 15).Value = [double]"0.005104969278514056"
$ws.Cells.Item(7, 16).Value = [double]"103.276926599096"
$ws.Cells.Item(7, 17).Value = [double]"153.2508354186922"

$ws.Cells.Item(8, 1).Value = "model_22_3_18"
$ws.Cells.Item(8, 2).Value = [double]"0.999974724567298"
$ws.Cells.Item(8, 3).Value = [double]"0.9990496085086725"
$ws.Cells.Item(8, 4).Value = [double]"0.9999056547386764"
$ws.Cells.Item(8, 5).Value = [double]"0.999978456303359"
$ws.Cells.Item(8, 6).Value = [double]"0.9999261340664485"
$ws.Cells.Item(8, 7).Value = [double]"2.359350696410027e-05"
$ws.Cells.Item(8, 8).Value = [double]"0.0008871487397891924"
$ws.Cells.Item(8, 9).Value = [double]"0.0001316446199193724"
$ws.Cells.Item(8, 10).Value = [double]"1.176502249724065e-05"
$ws.Cells.Item(8, 11).Value = [double]"7.170482120830655e-05"
$ws.Cells.Item(8, 12).Value = [double]"0.0003403224526816273"
$ws.Cells.Item(8, 13).Value = [double]"0.004857314789479911"
$ws.Cells.Item(8, 14).Value = [double]"1.000035682963815"
$ws.Cells.Item(8, 15).Value = [double]"0.005064100736095836"
$ws.Cells.Item(8, 16).Value = [double]"103.3090780248547"
$ws.Cells.Item(8, 17).Value = [double]"153.2829868444509"

$ws.Cells.Item(9, 1).Value = "model_22_3_17"
$ws.Cells.Item(9, 2).Value = [double]"0.9999751589247136"
$ws.Cells.Item(9, 3).Value = [double]"0.9990476822645158"
$ws.Cells.Item(9, 4).Value = [double]"0.9999082077594839"
$ws.Cells.Item(9, 5).Value = [double]"0.9999782903589507"
$ws.Cells.Item(9, 6).Value = [double]"0.9999279222508504"
$ws.Cells.Item(9, 7).Value = [double]"2.318805338266505e-05"
$ws.Cells.Item(9, 8).Value = [double]"0.0008889468041571923"
$ws.Cells.Item(9, 9).Value = [double]"0.0001280822634308325"
$ws.Cells.Item(9, 10).Value = [double]"1.185564481381522e-05"
$ws.Cells.Item(9, 11).Value = [double]"6.996895412232384e-05"
$ws.Cells.Item(9, 12).Value = [double]"0.0003453184506954858"
$ws.Cells.Item(9, 13).Value = [double]"0.004815397531114649"
$ws.Cells.Item(9, 14).Value = [double]"1.000035069753346"
$ws.Cells.Item(9, 15).Value = [double]"0.005020398973261277"
$ws.Cells.Item(9, 16).Value = [double]"103.3437467046434"
$ws.Cells.Item(9, 17).Value = [double]"153.3176555242397"

$ws.Cells.Item(10, 1).Value = "model_22_3_16"
$ws.Cells.Item(10, 2).Value = [double]"0.9999756029564093"
$ws.Cells.Item(10, 3).Value = [double]"0.9990455058910916"
$ws.Cells.Item(10, 4).Value = [double]"0.999910974387391"
$ws.Cells.Item(10, 5).Value = [double]"0.9999780725647875"
$ws.Cells.Item(10, 6).Value = [double]"0.9999298493708219"
$ws.Cells.Item(10, 7).Value = [double]"2.277356928551199e-05"
$ws.Cells.Item(10, 8).Value = [double]"0.0008909783532168728"
$ws.Cells.Item(10, 9).Value = [double]"0.0001242218503672732"
$ws.Cells.Item(10, 10).Value = [double]"1.197458230504662e-05"
$ws.Cells.Item(10, 11).Value = [double]"6.809821633615991e-05"
$ws.Cells.Item(10, 12).Value = [double]"0.0003508087589167091"
$ws.Cells.Item(10, 13).Value = [double]"0.00477216609995"
$ws.Cells.Item(10, 14).Value = [double]"1.000034442885069"
$ws.Cells.Item(10, 15).Value = [double]"0.004975327090570549"
$ws.Cells.Item(10, 16).Value = [double]"103.3798198726159"
$ws.Cells.Item(10, 17).Value = [double]"153.3537286922122"

$ws.Cells.Item(11, 1).Value = "model_22_3_15"
$ws.Cells.Item(11, 2).Value = [double]"0.999976059322864"
$ws.Cells.Item(11, 3).Value = [double]"0.9990430411789265"
$ws.Cells.Item(11, 4).Value = [double]"0.9999139989608272"
$ws.Cells.Item(11, 5).Value = [double]"0.9999778040777337"
$ws.Cells.Item(11, 6).Value = [double]"0.9999319476178838"
$ws.Cells.Item(11, 7).Value = [double]"2.234757123213243e-05"
$ws.Cells.Item(11, 8).Value = [double]"0.0008932790538346913"
$ws.Cells.Item(11, 9).Value = [double]"0.0001200015131203554"
$ws.Cells.Item(11, 10).Value = [double]"1.212120320673745e-05"
$ws.Cells.Item(11, 11).Value = [double]"6.606135816354642e-05"
$ws.Cells.Item(11, 12).Value = [double]"0.0003567235399110237"
$ws.Cells.Item(11, 13).Value = [double]"0.00472732178216508"
$ws.Cells.Item(11, 14).Value = [double]"1.000033798603015"
$ws.Cells.Item(11, 15).Value = [double]"0.004928573657336991"
$ws.Cells.Item(11, 16).Value = [double]"103.4175858250258"
$ws.Cells.Item(11, 17).Value = [double]"153.3914946446221"

$ws.Cells.Item(12, 1).Value = "model_22_3_14"
$ws.Cells.Item(12, 2).Value = [double]"0.9999765174009664"
$ws.Cells.Item(12, 3).Value = [double]"0.9990402266960854"
$ws.Cells.Item(12, 4).Value = [double]"0.9999172841917062"
$ws.Cells.Item(12, 5).Value = [double]"0.9999774725938742"
$ws.Cells.Item(12, 6).Value = [double]"0.9999342154802526"
$ws.Cells.Item(12, 7).Value = [double]"2.191997543096203e-05"
$ws.Cells.Item(12, 8).Value = [double]"0.0008959062500253067"
$ws.Cells.Item(12, 9).Value = [double]"0.0001154174676226117"
$ws.Cells.Item(12, 10).Value = [double]"1.230222669266381e-05"
$ws.Cells.Item(12, 11).Value = [double]"6.385984715763778e-05"
$ws.Cells.Item(12, 12).Value = [double]"0.0003631305446162032"
$ws.Cells.Item(12, 13).Value = [double]"0.004681877340443898"
$ws.Cells.Item(12, 14).Value = [double]"1.000033151904518"
$ws.Cells.Item(12, 15).Value = [double]"0.004881194551648775"
$ws.Cells.Item(12, 16).Value = [double]"103.4562244334707"
$ws.Cells.Item(12, 17).Value = [double]"153.4301332530669"

$ws.Cells.Item(13, 1).Value = "model_22_3_13"
$ws.Cells.Item(13, 2).Value = [double]"0.9999769649683165"
$ws.Cells.Item(13, 3).Value = [double]"0.9990370589383444"
$ws.Cells.Item(13, 4).Value = [double]"0.9999208372467289"
$ws.Cells.Item(13, 5).Value = [double]"0.9999770549791474"
$ws.Cells.Item(13, 6).Value = [double]"0.9999366516015089"
$ws.Cells.Item(13, 7).Value = [double]"2.150219095551268e-05"
$ws.Cells.Item(13, 8).Value = [double]"0.0008988632128280025"
$ws.Cells.Item(13, 9).Value = [double]"0.0001104597138208118"
$ws.Cells.Item(13, 10).Value = [double]"1.253028628423353e-05"
$ws.Cells.Item(13, 11).Value = [double]"6.149500005252266e-05"
$ws.Cells.Item(13, 12).Value = [double]"0.0003701064505580868"
$ws.Cells.Item(13, 13).Value = [double]"0.004637045498538124"
$ws.Cells.Item(13, 14).Value = [double]"1.00003252004473"
$ws.Cells.Item(13, 15).Value = [double]"0.004834454125418365"
$ws.Cells.Item(13, 16).Value = [double]"103.4947114462301"
$ws.Cells.Item(13, 17).Value = [double]"153.4686202658263"

$ws.Cells.Item(14, 1).Value = "model_22_3_12"
$ws.Cells.Item(14, 2).Value = [double]"0.9999774015715477"
$ws.Cells.Item(14, 3).Value = [double]"0.999033513684626"
$ws.Cells.Item(14, 4).Value = [double]"0.9999247066555863"
$ws.Cells.Item(14, 5).Value = [double]"0.9999765772572711"
$ws.Cells.Item(14, 6).Value = [double]"0.9999392981799828"
$ws.Cells.Item(14, 7).Value = [double]"2.109464100390959e-05"
$ws.Cells.Item(14, 8).Value = [double]"0.0009021725515554589"
$ws.Cells.Item(14, 9).Value = [double]"0.0001050605358312638"
$ws.Cells.Item(14, 10).Value = [double]"1.279117041733464e-05"
$ws.Cells.Item(14, 11).Value = [double]"5.892585312429922e-05"
$ws.Cells.Item(14, 12).Value = [double]"0.0003775316313907666"
$ws.Cells.Item(14, 13).Value = [double]"0.004592890266913591"
$ws.Cells.Item(14, 14).Value = [double]"1.000031903663697"
$ws.Cells.Item(14, 15).Value = [double]"0.004788419114169645"
$ws.Cells.Item(14, 16).Value = [double]"103.5329830611964"
$ws.Cells.Item(14, 17).Value = [double]"153.5068918807927"

$ws.Cells.Item(15, 1).Value = "model_22_3_11"
$ws.Cells.Item(15, 2).Value = [double]"0.9999777859909298"
$ws.Cells.Item(15, 3).Value = [double]"0.9990294380525899"
$ws.Cells.Item(15, 4).Value = [double]"0.9999288292901795"
$ws.Cells.Item(15, 5).Value = [double]"0.9999759669141352"
$ws.Cells.Item(15, 6).Value = [double]"0.9999420894486359"
$ws.Cells.Item(15, 7).Value = [double]"2.073580238474958e-05"
$ws.Cells.Item(15, 8).Value = [double]"0.0009059769751616023"
$ws.Cells.Item(15, 9).Value = [double]"9.93080194200797e-05"
$ws.Cells.Item(15, 10).Value = [double]"1.312447908040606e-05"
$ws.Cells.Item(15, 11).Value = [double]"5.621624925024287e-05"
$ws.Cells.Item(15, 12).Value = [double]"0.0003856135485441103"
$ws.Cells.Item(15, 13).Value = [double]"0.004553658132177863"
$ws.Cells.Item(15, 14).Value = [double]"1.000031360953981"
$ws.Cells.Item(15, 15).Value = [double]"0.004747516786236502"
$ws.Cells.Item(15, 16).Value = [double]"103.5672975357925"
$ws.Cells.Item(15, 17).Value = [double]"153.5412063553888"

$ws.Cells.Item(16, 1).Value = "model_22_3_10"
$ws.Cells.Item(16, 2).Value = [double]"0.9999781236303374"
$ws.Cells.Item(16, 3).Value = [double]"0.9990248399413394"
$ws.Cells.Item(16, 4).Value = [double]"0.9999332981260873"
$ws.Cells.Item(16, 5).Value = [double]"0.999975242260785"
$ws.Cells.Item(16, 6).Value = [double]"0.9999450973797291"
$ws.Cells.Item(16, 7).Value = [double]"2.042063081838798e-05"
$ws.Cells.Item(16, 8).Value = [double]"0.0009102691101802425"
$ws.Cells.Item(16, 9).Value = [double]"9.307243115292618e-05"
$ws.Cells.Item(16, 10).Value = [double]"1.352021260323241e-05"
$ws.Cells.Item(16, 11).Value = [double]"5.329632187807929e-05"
$ws.Cells.Item(16, 12).Value = [double]"0.0003942677709199689"
$ws.Cells.Item(16, 13).Value = [double]"0.004518919209101661"
$ws.Cells.Item(16, 14).Value = [double]"1.000030884286582"
$ws.Cells.Item(16, 15).Value = [double]"0.004711298955285463"
$ws.Cells.Item(16, 16).Value = [double]"103.5979297070427"
$ws.Cells.Item(16, 17).Value = [double]"153.571838526639"

$ws.Cells.Item(17, 1).Value = "model_22_3_9"
$ws.Cells.Item(17, 2).Value = [double]"0.999978372725472"
$ws.Cells.Item(17, 3).Value = [double]"0.9990196264008927"
$ws.Cells.Item(17, 4).Value = [double]"0.9999380526974727"
$ws.Cells.Item(17, 5).Value = [double]"0.9999743954083631"
$ws.Cells.Item(17, 6).Value = [double]"0.9999482762973683"
$ws.Cells.Item(17, 7).Value = [double]"2.018811144424163e-05"
$ws.Cells.Item(17, 8).Value = [double]"0.000915135721339236"
$ws.Cells.Item(17, 9).Value = [double]"8.643814200975663e-05"
$ws.Cells.Item(17, 10).Value = [double]"1.398267909457863e-05"
$ws.Cells.Item(17, 11).Value = [double]"5.021041055216763e-05"
$ws.Cells.Item(17, 12).Value = [double]"0.0004034412796105347"
$ws.Cells.Item(17, 13).Value = [double]"0.004493118231723001"
$ws.Cells.Item(17, 14).Value = [double]"1.000030532622863"
$ws.Cells.Item(17, 15).Value = [double]"0.004684399576884407"
$ws.Cells.Item(17, 16).Value = [double]"103.6208333383609"
$ws.Cells.Item(17, 17).Value = [double]"153.5947421579571"

$ws.Cells.Item(18, 1).Value = "model_22_3_8"
$ws.Cells.Item(18, 2).Value = [double]"0.9999785118223954"
$ws.Cells.Item(18, 3).Value = [double]"0.9990138005881417"
$ws.Cells.Item(18, 4).Value = [double]"0.9999431098754612"
$ws.Cells.Item(18, 5).Value = [double]"0.9999734546581784"
$ws.Cells.Item(18, 6).Value = [double]"0.9999516462874486"
$ws.Cells.Item(18, 7).Value = [double]"2.005827057187077e-05"
$ws.Cells.Item(18, 8).Value = [double]"0.0009205738618186328"
$ws.Cells.Item(18, 9).Value = [double]"7.938161087267512e-05"
$ws.Cells.Item(18, 10).Value = [double]"1.449642319667256e-05"
$ws.Cells.Item(18, 11).Value = [double]"4.693901703467383e-05"
$ws.Cells.Item(18, 12).Value = [double]"0.0004128652407990863"
$ws.Cells.Item(18, 13).Value = [double]"0.004478646064590365"
$ws.Cells.Item(18, 14).Value = [double]"1.000030336250736"
$ws.Cells.Item(18, 15).Value = [double]"0.00466931129963542"
$ws.Cells.Item(18, 16).Value = [double]"103.6337379838304"
$ws.Cells.Item(18, 17).Value = [double]"153.6076468034266"

$ws.Cells.Item(19, 1).Value = "model_22_3_7"
$ws.Cells.Item(19, 2).Value = [double]"0.999978482356387"
$ws.Cells.Item(19, 3).Value = [double]"0.9990072895954343"
$ws.Cells.Item(19, 4).Value = [double]"0.9999483892292284"
$ws.Cells.Item(19, 5).Value = [double]"0.9999724255603002"
$ws.Cells.Item(19, 6).Value = [double]"0.9999551491692276"
$ws.Cells.Item(19, 7).Value = [double]"2.008577579732539e-05"
$ws.Cells.Item(19, 8).Value = [double]"0.0009266515877115366"
$ws.Cells.Item(19, 9).Value = [double]"7.201506685829263e-05"
$ws.Cells.Item(19, 10).Value = [double]"1.505841401429894e-05"
$ws.Cells.Item(19, 11).Value = [double]"4.353861986099022e-05"
$ws.Cells.Item(19, 12).Value = [double]"0.0004224908007473032"
$ws.Cells.Item(19, 13).Value = [double]"0.004481715720271132"
$ws.Cells.Item(19, 14).Value = [double]"1.000030377849807"
$ws.Cells.Item(19, 15).Value = [double]"0.004672511636913578"
$ws.Cells.Item(19, 16).Value = [double]"103.6309973303839"
$ws.Cells.Item(19, 17).Value = [double]"153.6049061499802"

$ws.Cells.Item(20, 1).Value = "model_22_3_6"
$ws.Cells.Item(20, 2).Value = [double]"0.9999782717594072"
$ws.Cells.Item(20, 3).Value = [double]"0.999000107603992"
$ws.Cells.Item(20, 4).Value = [double]"0.9999540137121851"
$ws.Cells.Item(20, 5).Value = [double]"0.9999713655589163"
$ws.Cells.Item(20, 6).Value = [double]"0.999958895275659"
$ws.Cells.Item(20, 7).Value = [double]"2.028235883382174e-05"
$ws.Cells.Item(20, 8).Value = [double]"0.0009333556614698276"
$ws.Cells.Item(20, 9).Value = [double]"6.416694697720249e-05"
$ws.Cells.Item(20, 10).Value = [double]"1.563728125030418e-05"
$ws.Cells.Item(20, 11).Value = [double]"3.990211411375334e-05"
$ws.Cells.Item(20, 12).Value = [double]"0.0004320518792775426"
$ws.Cells.Item(20, 13).Value = [double]"0.004503593990783554"
$ws.Cells.Item(20, 14).Value = [double]"1.00003067516319"
$ws.Cells.Item(20, 15).Value = [double]"0.004695321310695979"
$ws.Cells.Item(20, 16).Value = [double]"103.6115181454007"
$ws.Cells.Item(20, 17).Value = [double]"153.5854269649969"

$ws.Cells.Item(21, 1).Value = "model_22_3_5"
$ws.Cells.Item(21, 2).Value = [double]"0.9999778251368089"
$ws.Cells.Item(21, 3).Value = [double]"0.9989923572199363"
$ws.Cells.Item(21, 4).Value = [double]"0.9999599122255013"
$ws.Cells.Item(21, 5).Value = [double]"0.9999704608984585"
$ws.Cells.Item(21, 6).Value = [double]"0.9999628800864401"
$ws.Cells.Item(21, 7).Value = [double]"2.069926142493499e-05"
$ws.Cells.Item(21, 8).Value = [double]"0.000940590304783288"
$ws.Cells.Item(21, 9).Value = [double]"5.593645895151942e-05"
$ws.Cells.Item(21, 10).Value = [double]"1.613131673618223e-05"
$ws.Cells.Item(21, 11).Value = [double]"3.603388784385083e-05"
$ws.Cells.Item(21, 12).Value = [double]"0.0004412061833730867"
$ws.Cells.Item(21, 13).Value = [double]"0.004549644098710909"
$ws.Cells.Item(21, 14).Value = [double]"1.000031305689211"
$ws.Cells.Item(21, 15).Value = [double]"0.004743331867054666"
$ws.Cells.Item(21, 16).Value = [double]"103.5708250765687"
$ws.Cells.Item(21, 17).Value = [double]"153.5447338961649"

$ws.Cells.Item(22, 1).Value = "model_22_3_4"
$ws.Cells.Item(22, 2).Value = [double]"0.9999770769401084"
$ws.Cells.Item(22, 3).Value = [double]"0.9989843480174928"
$ws.Cells.Item(22, 4).Value = [double]"0.999966041794364"
$ws.Cells.Item(22, 5).Value = [double]"0.9999700155026697"
$ws.Cells.Item(22, 6).Value = [double]"0.9999671601391315"
$ws.Cells.Item(22, 7).Value = [double]"2.139767020280001e-05"
$ws.Cells.Item(22, 8).Value = [double]"0.0009480665437010052"
$ws.Cells.Item(22, 9).Value = [double]"4.738356766820788e-05"
$ws.Cells.Item(22, 10).Value = [double]"1.637454757827333e-05"
$ws.Cells.Item(22, 11).Value = [double]"3.18790576232406e-05"
$ws.Cells.Item(22, 12).Value = [double]"0.0004502862003883709"
$ws.Cells.Item(22, 13).Value = [double]"0.004625761580842663"
$ws.Cells.Item(22, 14).Value = [double]"1.000032361966906"
$ws.Cells.Item(22, 15).Value = [double]"0.004822689827985679"
$ws.Cells.Item(22, 16).Value = [double]"103.5044570217823"
$ws.Cells.Item(22, 17).Value = [double]"153.4783658413786"

$ws.Cells.Item(23, 1).Value = "model_22_3_3"
$ws.Cells.Item(23, 2).Value = [double]"0.9999759148669479"
$ws.Cells.Item(23, 3).Value = [double]"0.9989762871816444"
$ws.Cells.Item(23, 4).Value = [double]"0.9999723259748917"
$ws.Cells.Item(23, 5).Value = [double]"0.999970284464537"
$ws.Cells.Item(23, 6).Value = [double]"0.9999717522458379"
$ws.Cells.Item(23, 7).Value = [double]"2.248241448898537e-05"
$ws.Cells.Item(23, 8).Value = [double]"0.0009555909801356036"
$ws.Cells.Item(23, 9).Value = [double]"3.861493906434759e-05"
$ws.Cells.Item(23, 10).Value = [double]"1.622766738070652e-05"
$ws.Cells.Item(23, 11).Value = [double]"2.742130322252706e-05"
$ws.Cells.Item(23, 12).Value = [double]"0.0004603057517614868"
$ws.Cells.Item(23, 13).Value = [double]"0.004741562452292005"
$ws.Cells.Item(23, 14).Value = [double]"1.000034002540779"
$ws.Cells.Item(23, 15).Value = [double]"0.004943420582273469"
$ws.Cells.Item(23, 16).Value = [double]"103.4055542652255"
$ws.Cells.Item(23, 17).Value = [double]"153.3794630848217"

$ws.Cells.Item(24, 1).Value = "model_22_3_2"
$ws.Cells.Item(24, 2).Value = [double]"0.999974200771291"
$ws.Cells.Item(24, 3).Value = [double]"0.998968899859716"
$ws.Cells.Item(24, 4).Value = [double]"0.9999784548562639"
$ws.Cells.Item(24, 5).Value = [double]"0.9999718230121192"
$ws.Cells.Item(24, 6).Value = [double]"0.999976589847405"
$ws.Cells.Item(24, 7).Value = [double]"2.408244754462745e-05"
$ws.Cells.Item(24, 8).Value = [double]"0.0009624867208898475"
$ws.Cells.Item(24, 9).Value = [double]"3.006300707054296e-05"
$ws.Cells.Item(24, 10).Value = [double]"1.538746584895418e-05"
$ws.Cells.Item(24, 11).Value = [double]"2.272523645974857e-05"
$ws.Cells.Item(24, 12).Value = [double]"0.0004704998785817867"
$ws.Cells.Item(24, 13).Value = [double]"0.004907387038397059"
$ws.Cells.Item(24, 14).Value = [double]"1.00003642244053"
$ws.Cells.Item(24, 15).Value = [double]"0.005116304664313231"
$ws.Cells.Item(24, 16).Value = [double]"103.2680546009402"
$ws.Cells.Item(24, 17).Value = [double]"153.2419634205365"

$ws.Cells.Item(25, 1).Value = "model_22_3_1"
$ws.Cells.Item(25, 2).Value = [double]"0.999971761450658"
$ws.Cells.Item(25, 3).Value = [double]"0.998962861485615"
$ws.Cells.Item(25, 4).Value = [double]"0.9999841267330543"
$ws.Cells.Item(25, 5).Value = [double]"0.9999751319155743"
$ws.Cells.Item(25, 6).Value = [double]"0.9999815969640107"
$ws.Cells.Item(25, 7).Value = [double]"2.635944628167136e-05"
$ws.Cells.Item(25, 8).Value = [double]"0.0009681232780592204"
$ws.Cells.Item(25, 9).Value = [double]"2.214875622402684e-05"
$ws.Cells.Item(25, 10).Value = [double]"1.358047217280066e-05"
$ws.Cells.Item(25, 11).Value = [double]"1.786461419841374e-05"
$ws.Cells.Item(25, 12).Value = [double]"0.0004801879210398119"
$ws.Cells.Item(25, 13).Value = [double]"0.005134145136405024"
$ws.Cells.Item(25, 14).Value = [double]"1.000039866187306"
$ws.Cells.Item(25, 15).Value = [double]"0.005352716324007369"
$ws.Cells.Item(25, 16).Value = [double]"103.087367709114"
$ws.Cells.Item(25, 17).Value = [double]"153.0612765287102"

$ws.Cells.Item(26, 1).Value = "model_22_3_0"
$ws.Cells.Item(26, 2).Value = [double]"0.9999684454911321"
$ws.Cells.Item(26, 3).Value = [double]"0.9989584569617904"
$ws.Cells.Item(26, 4).Value = [double]"0.9999893746777766"
$ws.Cells.Item(26, 5).Value = [double]"0.9999803782924166"
$ws.Cells.Item(26, 6).Value = [double]"0.9999868452261689"
$ws.Cells.Item(26, 7).Value = [double]"2.945474894539345e-05"
$ws.Cells.Item(26, 8).Value = [double]"0.0009722347077131216"
$ws.Cells.Item(26, 9).Value = [double]"1.482603880674469e-05"
$ws.Cells.Item(26, 10).Value = [double]"1.071542340204817e-05"
$ws.Cells.Item(26, 11).Value = [double]"1.276990163449178e-05"
$ws.Cells.Item(26, 12).Value = [double]"0.0004880703027942576"
$ws.Cells.Item(26, 13).Value = [double]"0.005427222949667118"
$ws.Cells.Item(26, 14).Value = [double]"1.000044547541931"
$ws.Cells.Item(26, 15).Value = [double]"0.005658271066534741"
$ws.Cells.Item(26, 16).Value = [double]"102.8653108127086"
$ws.Cells.Item(26, 17).Value = [double]"152.8392196323048"
